$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to Text format, assign the literal string value, then restore
# the default ("Normal") style so no stray number-format/style artifact is left
# behind (the source data are text strings like "1.000" / "0.07699" that Excel
# would otherwise auto-coerce to numbers and reformat).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '28.934.24'
$ws.Range("E2").Value = '  -1.53%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.832.84'
$ws.Range("E3").Value = '  -1.89%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.9999'
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
Set-TextValue $ws.Range("D5") '245.43'
$ws.Range("E5").Value = '  +0.68%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.6902'
$ws.Range("E6").Value = '  -1.96%  '

# Row 7
$ws.Range("E7").Value = '  -0.09%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.07699'
$ws.Range("E8").Value = '  -2.84%  '

# Row 9
$ws.Range("E9").Value = '  -2.59%  '

# Row 10
Set-TextValue $ws.Range("D10") '23.49'
$ws.Range("E10").Value = '  -3.85%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.07819'
$ws.Range("E11").Value = '  -0.40%  '

# Row 12
Set-TextValue $ws.Range("D12") '1.839.99'
$ws.Range("E12").Value = '  -1.35%  '

# Row 13
Set-TextValue $ws.Range("D13") '5.082'
$ws.Range("E13").Value = '  -1.83%  '

# Row 14
Set-TextValue $ws.Range("D14") '90.45'
$ws.Range("E14").Value = '  -3.57%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.6805'
$ws.Range("E15").Value = '  -2.84%  '

# Row 16
Set-TextValue $ws.Range("D16") '6.440'
$ws.Range("E16").Value = '  -1.22%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.000008327'
$ws.Range("E17").Value = '  -0.83%  '

# Row 18
Set-TextValue $ws.Range("D18") '28.941.78'
$ws.Range("E18").Value = '  -1.49%  '

# Row 19
Set-TextValue $ws.Range("D19") '243.18'
$ws.Range("E19").Value = '  -4.24%  '

# Row 20
Set-TextValue $ws.Range("D20") '2.082.75'
$ws.Range("E20").Value = '  -1.46%  '

# Row 21
Set-TextValue $ws.Range("D21") '12.74'
$ws.Range("E21").Value = '  -2.78%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.9998'
$ws.Range("E22").Value = '  -0.02%  '

# Row 23
Set-TextValue $ws.Range("D23") '7.471'
$ws.Range("E23").Value = '  -2.29%  '

# Row 24
Set-TextValue $ws.Range("D24") '1.000'
$ws.Range("E24").Value = '  -0.12%  '

# Row 25
Set-TextValue $ws.Range("D25") '163.25'
$ws.Range("E25").Value = '  +1.13%  '

# Row 26
$ws.Range("E26").Value = '  -5.45%  '

# Row 27
Set-TextValue $ws.Range("D27") '8.801'
$ws.Range("E27").Value = '  -2.22%  '

# Row 28
Set-TextValue $ws.Range("D28") '18.20'
$ws.Range("E28").Value = '  -3.27%  '

# Row 29
Set-TextValue $ws.Range("D29") '1.551'
$ws.Range("E29").Value = '  +3.40%  '

# Row 30
Set-TextValue $ws.Range("D30") '4.212'
$ws.Range("E30").Value = '  -2.40%  '

# Row 31
Set-TextValue $ws.Range("D31") '4.156'
$ws.Range("E31").Value = '  -2.09%  '

# Row 32
Set-TextValue $ws.Range("D32") '1.175'
$ws.Range("E32").Value = '  -3.12%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.05114'
$ws.Range("E33").Value = '  -3.11%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.7712'
$ws.Range("E34").Value = '  +3.19%  '

# Row 35
Set-TextValue $ws.Range("D35") '1.842'
$ws.Range("E35").Value = '  -2.91%  '

# Row 36
$ws.Range("E36").Value = '  -2.63%  '

# Row 37
Set-TextValue $ws.Range("D37") '2.683'
$ws.Range("E37").Value = '  -0.94%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.01849'
$ws.Range("E38").Value = '  -1.61%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.242.19'
$ws.Range("E39").Value = '  -2.66%  '

# Row 40
Set-TextValue $ws.Range("D40") '2.696'
$ws.Range("E40").Value = '  -2.51%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.9408'
$ws.Range("E41").Value = '  +5.51%  '

# Row 42
Set-TextValue $ws.Range("D42") '108.17'
$ws.Range("E42").Value = '  -0.37%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.9996'
$ws.Range("E43").Value = '  -0.13%  '

# Row 44
Set-TextValue $ws.Range("D44") '5.672'
$ws.Range("E44").Value = '  -5.57%  '

# Row 45
Set-TextValue $ws.Range("D45") '9.637'
$ws.Range("E45").Value = '  +0.35%  '

# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range("D46") '1.983.21'
$ws.Range("E46").Value = '  -1.59%  '

# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D47") '0.5173'
$ws.Range("E47").Value = '  -0.15%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D48") '0.00000000122'
$ws.Range("E48").Value = '  -5.27%  '

# Row 49
Set-TextValue $ws.Range("D49") '64.59'

# Row 50
Set-TextValue $ws.Range("D50") '1.747'
$ws.Range("E50").Value = '  -2.68%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.4197'
$ws.Range("E51").Value = '  -2.37%  '
